$wb = $excel.ActiveWorkbook

# The "Slovakia" sheet is the last country sheet in the workbook; the new
# "Hungary" sheet is a sibling tab built the same way (same layout/styles),
# so we clone it and then tweak the handful of cells that differ.
$slovakia = $wb.Worksheets.Item("Slovakia")

# Select the whole sheet on Slovakia before copying away from it - this
# mirrors the "previous selection" left behind on Slovakia once it is no
# longer the active tab.
$slovakia.Range("A1:XFD1048576").Select()

# Worksheet.Copy duplicates the sheet (styles, merges, column widths, etc.)
# and places the copy immediately after the source sheet; the copy becomes
# the new ActiveSheet/ActiveWorkbook tab.
$slovakia.Copy($null, $slovakia)
$hungary = $wb.ActiveSheet
$hungary.Name = "Hungary"

# Give row 4's data cell the same bordered style used by every other
# country sheet (Slovakia's B4 was the one outlier without it) by copying
# the formatting already used one row up.
$hungary.Range("B3").Copy()
$hungary.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Market name + ticket reference for the new Hungary market.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3593/T3618/T3592/T3617/T3594/T3619"

# Column B needs to widen to fit the longer reference text.
$hungary.Columns.Item(2).ColumnWidth = 44

# Leave the new sheet's selection on the cell that was just filled in.
$hungary.Range("B7").Select()
